$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values of columns D, K, L, M, O, P between row 2 and row 3.
$cols = @("D", "K", "L", "M", "O", "P")

foreach ($col in $cols) {
    $addr2 = "$col" + "2"
    $addr3 = "$col" + "3"
    $v2 = $ws.Range($addr2).Value2
    $v3 = $ws.Range($addr3).Value2
    $ws.Range($addr2).Value2 = $v3
    $ws.Range($addr3).Value2 = $v2
}
